# TurtleBomb_Stage_Table.xlsx - apply the commit's data/view changes to
# the "Table_Stage" sheet.
#
# Per the diff:
#   - C8: 7  -> 17
#   - D8: 8  -> 0
#   - the active cell/selection moves from H10 to C8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_Stage")

# --- Data edits (row 8 = stage 5 pattern IDs) ---
$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 0

# --- Selection / active cell moves to C8 ---
$ws.Activate()
$ws.Range("C8").Select()
